$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# RF002 row: "Pode chamar RF008" -> "Pode chamar RF007"
$ws.Range("E11").Value = "Pode chamar RF007"

# RF005 row: "Pode chamar RF006" -> "Pode chamar RF008"
$ws.Range("E14").Value = "Pode chamar RF008"

# Update the active selection shown in the saved view
$ws.Range("C1").Select()
